$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (PM2.5 -> NOx); this also updates the _FilterDatabase
# defined name automatically since it references the sheet by name.
$ws.Name = "2019_NOx"

# Update the pollutant value for the data row (PM2,5 -> NOx)
$ws.Range("C2").Value = "NOx"

# Match the selection left behind after the edit
$ws.Range("C3").Select() | Out-Null
